# All Country Files Saved And Formatted
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Rename header for column L
# ---------------------------------------------------------------------
$ws.Range("L1").Value = "MasterSheet RowNo."

# ---------------------------------------------------------------------
# 2. Append new data row (row 5) for serial date 43921 (2020-04-25)
# ---------------------------------------------------------------------
$ws.Cells.Item(5, 1).Value = 71
$ws.Cells.Item(5, 2).Value = 261
$ws.Cells.Item(5, 3).Value = "UNITED STATES VIRGIN ISLANDS"
$ws.Cells.Item(5, 4).Value = 43921
$ws.Cells.Item(5, 5).Value = "United States Virgin Islands"
$ws.Cells.Item(5, 6).Value = 30
$ws.Cells.Item(5, 7).Value = 8
$ws.Cells.Item(5, 8).Value = 0
$ws.Cells.Item(5, 9).Value = 0
$ws.Cells.Item(5, 10).Value = "Imported cases only"
$ws.Cells.Item(5, 11).Value = 0
$ws.Cells.Item(5, 12).Value = 5375

# ---------------------------------------------------------------------
# 3. Fix up existing rows 2-4: uppercase WorldRegion, fill in the
#    TotalConfirmedNewCases (G) / TotalNewDeaths (I) columns
# ---------------------------------------------------------------------
$ws.Cells.Item(2, 3).Value = "UNITED STATES VIRGIN ISLANDS"
$ws.Cells.Item(2, 7).Value = 17
$ws.Cells.Item(2, 9).Value = 0

$ws.Cells.Item(3, 3).Value = "UNITED STATES VIRGIN ISLANDS"
$ws.Cells.Item(3, 7).Value = 0
$ws.Cells.Item(3, 9).Value = 0

$ws.Cells.Item(4, 3).Value = "UNITED STATES VIRGIN ISLANDS"
$ws.Cells.Item(4, 7).Value = 5
$ws.Cells.Item(4, 9).Value = 0

# ---------------------------------------------------------------------
# 4. Column widths: A..O (1..15) -> 27 characters wide
# ---------------------------------------------------------------------
$ws.Range("A1:O1").EntireColumn.ColumnWidth = 26.14

# ---------------------------------------------------------------------
# 5. Formatting: whole used range (A1:O5) centered horizontally/vertically,
#    and the date column (D) additionally gets a custom date number format.
#    Build each target style once on a scratch cell, then fan it out with
#    a single format-only paste so the style table stays minimal/exact.
# ---------------------------------------------------------------------
$scratchAlign = $ws.Range("Z1")
$scratchAlign.HorizontalAlignment = -4108
$scratchAlign.VerticalAlignment = -4108
$scratchAlign.Copy()
$ws.Range("A1:O5").PasteSpecial(-4122)

$scratchDate = $ws.Range("Z2")
$scratchDate.NumberFormat = "yyyy-mm-dd;"
$scratchDate.HorizontalAlignment = -4108
$scratchDate.VerticalAlignment = -4108
$scratchDate.Copy()
$ws.Range("D1:D5").PasteSpecial(-4122)

# Remove the scratch cells used to build the styles above
$ws.Range("Z1:Z2").Clear()
